$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cryptocurrency price / volume data (scraped refresh) ---

$ws.Range("D2").Value = "68.628.10"
$ws.Range("E2").Value = "  -1.60%  "
$ws.Range("D3").Value = "2.456.68"
$ws.Range("E3").Value = "  -2.13%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "565.07"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.69"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.85%  "
$ws.Range("E8").Value = "  -1.35%  "
$ws.Range("E9").Value = "  -6.14%  "
$ws.Range("E10").Value = "  -1.94%  "
$ws.Range("E11").Value = "  -3.92%  "
$ws.Range("E12").Value = "  -2.48%  "
$ws.Range("D13").Value = "2.907.14"
$ws.Range("E13").Value = "  -2.21%  "
$ws.Range("D14").Value = "68.515.84"
$ws.Range("E14").Value = "  -1.68%  "
$ws.Range("E15").Value = "  -3.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.67"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -4.85%  "
$ws.Range("D17").Value = "2.469.47"
$ws.Range("E17").Value = "  -1.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.02"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "344.61"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.20"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.13%  "
$ws.Range("E21").Value = "  -2.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.89"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.86%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.21"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.89%  "
$ws.Range("E25").Value = "  -4.65%  "
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "2.582.87"
$ws.Range("E26").Value = "  -2.67%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.02"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.25"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -6.76%  "
$ws.Range("D29").Value = "0.0₃0843"
$ws.Range("E29").Value = "  -5.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.33"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -6.88%  "
$ws.Range("E31").Value = "  -3.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "436.82"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("E34").Value = "  -2.66%  "
$ws.Range("E35").Value = "  +101.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.07"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  -5.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.92"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.14%  "
$ws.Range("E41").Value = "  -3.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.52"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.83%  "
$ws.Range("E43").Value = "  -3.64%  "
$ws.Range("E44").Value = "  +1.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.10"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "135.48"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.64%  "
$ws.Range("E47").Value = "  -2.36%  "
$ws.Range("E48").Value = "  -5.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0720"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.565"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.55%  "
